# Append the 2025-12-28 bitcoin buy entry as a new row (row 63), matching
# the pattern already used by the other "manually logged" rows (e.g. rows
# 61-62) where the Date column is stored as a literal text string such as
# "12/21/2025" rather than a real Excel date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so Excel treats the value as literal text
# instead of auto-converting "12/28/2025" into a date. Then strip the
# resulting "quote prefix" formatting so the cell ends up with no special
# number format/style applied (consistent with the other text-date cells
# in this sheet).
$ws.Cells.Item(63, 1).Value = "'12/28/2025"
$ws.Cells.Item(63, 1).ClearFormats()

$ws.Cells.Item(63, 2).Value = 0.0005624099999999993
$ws.Cells.Item(63, 3).Value = 88014.08225316061
$ws.Cells.Item(63, 4).Value = 50
